# Auto-generated edit script: restores/updates cached numeric values
# in the "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns
# (H, I, J, K, L, M, N) for specific rows across all 8 craft-job sheets,
# per the scheduled-runner price refresh (Typhon_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17, J17, L17, N17
$ws.Cells.Item(17, 8).Value = 3595064.5
$ws.Cells.Item(17, 10).Value = 3728144.8
$ws.Cells.Item(17, 12).Value = 11184434.4
$ws.Cells.Item(17, 14).Value = -11184770.4
# Row 41: H41, I41, J41, K41, L41, M41, N41
$ws.Cells.Item(41, 8).Value = 674.1875
$ws.Cells.Item(41, 9).Value = 127
$ws.Cells.Item(41, 10).Value = 800.46155
$ws.Cells.Item(41, 11).Value = 127
$ws.Cells.Item(41, 12).Value = 800.46155
$ws.Cells.Item(41, 13).Value = 313
$ws.Cells.Item(41, 14).Value = -1680.46155
# Row 51: H51, J51, L51, N51
$ws.Cells.Item(51, 8).Value = 6416.6665
$ws.Cells.Item(51, 10).Value = 5001
$ws.Cells.Item(51, 12).Value = 5001
$ws.Cells.Item(51, 14).Value = -5969
# Row 53: H53, J53, L53, N53
$ws.Cells.Item(53, 8).Value = 1322.25
$ws.Cells.Item(53, 10).Value = 1679.6666
$ws.Cells.Item(53, 12).Value = 1679.6666
$ws.Cells.Item(53, 14).Value = -2953.6666
# Row 98: H98, I98, J98, K98, L98, M98, N98
$ws.Cells.Item(98, 8).Value = 594.7059
$ws.Cells.Item(98, 9).Value = 594.7059
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 594.7059
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 903.2941
$ws.Cells.Item(98, 14).ClearContents()
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Cells.Item(122, 8).Value = 594.7059
$ws.Cells.Item(122, 9).Value = 594.7059
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 1784.1177
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 665.8822999999998
$ws.Cells.Item(122, 14).ClearContents()
# Row 129: H129, J129, L129, N129
$ws.Cells.Item(129, 8).Value = 345777.1
$ws.Cells.Item(129, 10).Value = 401062.38
$ws.Cells.Item(129, 12).Value = 1203187.14
$ws.Cells.Item(129, 14).Value = -1213187.14
# Row 138: H138, I138, K138, M138
$ws.Cells.Item(138, 8).Value = 3144.383
$ws.Cells.Item(138, 9).Value = 2309.8823
$ws.Cells.Item(138, 11).Value = 6929.646900000001
$ws.Cells.Item(138, 13).Value = -1789.646900000001
# Row 141: H141, I141, K141, M141
$ws.Cells.Item(141, 8).Value = 2463.8064
$ws.Cells.Item(141, 9).Value = 2138.5
$ws.Cells.Item(141, 11).Value = 6415.5
$ws.Cells.Item(141, 13).Value = -1235.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, J32, L32, N32
$ws.Cells.Item(32, 8).Value = 10197.424
$ws.Cells.Item(32, 10).Value = 25808.889
$ws.Cells.Item(32, 12).Value = 25808.889
$ws.Cells.Item(32, 14).Value = -26382.889
# Row 45: H45, I45, K45, M45
$ws.Cells.Item(45, 8).Value = 2383.611
$ws.Cells.Item(45, 9).Value = 2001.1852
$ws.Cells.Item(45, 11).Value = 2001.1852
$ws.Cells.Item(45, 13).Value = -1624.1852
# Row 61: H61, I61, J61, K61, L61, M61, N61
$ws.Cells.Item(61, 8).Value = 10104435
$ws.Cells.Item(61, 9).Value = 12349129
$ws.Cells.Item(61, 10).Value = 3316.5
$ws.Cells.Item(61, 11).Value = 12349129
$ws.Cells.Item(61, 12).Value = 3316.5
$ws.Cells.Item(61, 13).Value = -12348917
$ws.Cells.Item(61, 14).Value = -3740.5
# Row 122: H122, I122, K122, M122
$ws.Cells.Item(122, 8).Value = 3417.8462
$ws.Cells.Item(122, 9).Value = 3411.0833
$ws.Cells.Item(122, 11).Value = 10233.2499
$ws.Cells.Item(122, 13).Value = -7783.249899999999
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Cells.Item(132, 8).Value = 9627240
$ws.Cells.Item(132, 9).Value = 12197069
$ws.Cells.Item(132, 10).Value = 48786.363
$ws.Cells.Item(132, 11).Value = 36591207
$ws.Cells.Item(132, 12).Value = 146359.089
$ws.Cells.Item(132, 13).Value = -36588677
$ws.Cells.Item(132, 14).Value = -151419.089
# Row 134: H134, J134, L134, N134
$ws.Cells.Item(134, 8).Value = 45000
$ws.Cells.Item(134, 10).Value = 45000
$ws.Cells.Item(134, 12).Value = 45000
$ws.Cells.Item(134, 14).Value = -55140
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Cells.Item(136, 8).Value = 10104435
$ws.Cells.Item(136, 9).Value = 12349129
$ws.Cells.Item(136, 10).Value = 3316.5
$ws.Cells.Item(136, 11).Value = 37047387
$ws.Cells.Item(136, 12).Value = 9949.5
$ws.Cells.Item(136, 13).Value = -37044837
$ws.Cells.Item(136, 14).Value = -15049.5

$ws = $wb.Worksheets.Item("BSM")
# Row 99: H99, I99, K99, M99
$ws.Cells.Item(99, 8).Value = 1402.2222
$ws.Cells.Item(99, 9).Value = 1502.8572
$ws.Cells.Item(99, 11).Value = 1502.8572
$ws.Cells.Item(99, 13).Value = -4.857199999999921
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Cells.Item(107, 8).Value = 2416.795
$ws.Cells.Item(107, 9).Value = 2081.4443
$ws.Cells.Item(107, 10).Value = 3171.3333
$ws.Cells.Item(107, 11).Value = 2081.4443
$ws.Cells.Item(107, 12).Value = 3171.3333
$ws.Cells.Item(107, 13).Value = -161.4443000000001
$ws.Cells.Item(107, 14).Value = -7011.3333
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Cells.Item(134, 8).Value = 4880.5835
$ws.Cells.Item(134, 9).Value = 5262.077
$ws.Cells.Item(134, 10).Value = 3888.7
$ws.Cells.Item(134, 11).Value = 15786.231
$ws.Cells.Item(134, 12).Value = 11666.1
$ws.Cells.Item(134, 13).Value = -13251.231
$ws.Cells.Item(134, 14).Value = -16736.1

$ws = $wb.Worksheets.Item("CRP")
# Row 28: H28, J28, L28, N28
$ws.Cells.Item(28, 8).Value = 17000
$ws.Cells.Item(28, 10).Value = 17000
$ws.Cells.Item(28, 12).Value = 17000
$ws.Cells.Item(28, 14).Value = -17490
# Row 58: H58, J58, L58, N58
$ws.Cells.Item(58, 8).Value = 14299.667
$ws.Cells.Item(58, 10).Value = 22231.291
$ws.Cells.Item(58, 12).Value = 22231.291
$ws.Cells.Item(58, 14).Value = -22637.291
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Cells.Item(132, 8).Value = 25643232
$ws.Cells.Item(132, 9).Value = 40001596
$ws.Cells.Item(132, 10).Value = 3296.2856
$ws.Cells.Item(132, 11).Value = 120004788
$ws.Cells.Item(132, 12).Value = 9888.856800000001
$ws.Cells.Item(132, 13).Value = -120002258
$ws.Cells.Item(132, 14).Value = -14948.8568
# Row 136: H136, J136, L136, N136
$ws.Cells.Item(136, 8).Value = 14299.667
$ws.Cells.Item(136, 10).Value = 22231.291
$ws.Cells.Item(136, 12).Value = 66693.87300000001
$ws.Cells.Item(136, 14).Value = -71793.87300000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5, J5, L5, N5
$ws.Cells.Item(5, 8).Value = 1144.079
$ws.Cells.Item(5, 10).Value = 1883.75
$ws.Cells.Item(5, 12).Value = 5651.25
$ws.Cells.Item(5, 14).Value = -5875.25
# Row 103: H103, J103, L103, N103
$ws.Cells.Item(103, 8).Value = 3465.2778
$ws.Cells.Item(103, 10).Value = 9023.333000000001
$ws.Cells.Item(103, 12).Value = 27069.999
$ws.Cells.Item(103, 14).Value = -28827.999
# Row 124: H124, I124, J124, K124, L124, M124, N124
$ws.Cells.Item(124, 8).Value = 2366.6667
$ws.Cells.Item(124, 9).Value = 300
$ws.Cells.Item(124, 10).Value = 6500
$ws.Cells.Item(124, 11).Value = 900
$ws.Cells.Item(124, 12).Value = 19500
$ws.Cells.Item(124, 13).Value = 4010
$ws.Cells.Item(124, 14).Value = -29320
# Row 131: H131, J131, L131, N131
$ws.Cells.Item(131, 8).Value = 659.96
$ws.Cells.Item(131, 10).Value = 695.98865
$ws.Cells.Item(131, 12).Value = 2087.96595
$ws.Cells.Item(131, 14).Value = -12167.96595
# Row 135: H135, J135, L135, N135
$ws.Cells.Item(135, 8).Value = 1144.079
$ws.Cells.Item(135, 10).Value = 1883.75
$ws.Cells.Item(135, 12).Value = 16953.75
$ws.Cells.Item(135, 14).Value = -22023.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80, I80, J80, K80, L80, M80, N80
$ws.Cells.Item(80, 8).Value = 3836.818
$ws.Cells.Item(80, 9).Value = 3488
$ws.Cells.Item(80, 10).Value = 4036.1428
$ws.Cells.Item(80, 11).Value = 3488
$ws.Cells.Item(80, 12).Value = 4036.1428
$ws.Cells.Item(80, 13).Value = -2490
$ws.Cells.Item(80, 14).Value = -6032.1428
# Row 83: H83, I83, J83, K83, L83, M83, N83
$ws.Cells.Item(83, 8).Value = 3836.818
$ws.Cells.Item(83, 9).Value = 3488
$ws.Cells.Item(83, 10).Value = 4036.1428
$ws.Cells.Item(83, 11).Value = 17440
$ws.Cells.Item(83, 12).Value = 20180.714
$ws.Cells.Item(83, 13).Value = -12448
$ws.Cells.Item(83, 14).Value = -30164.714
# Row 97: H97, I97, J97, K97, L97, M97, N97
$ws.Cells.Item(97, 8).Value = 856.93335
$ws.Cells.Item(97, 9).Value = 417.35715
$ws.Cells.Item(97, 10).Value = 7011
$ws.Cells.Item(97, 11).Value = 417.35715
$ws.Cells.Item(97, 12).Value = 7011
$ws.Cells.Item(97, 13).Value = 78.64285000000001
$ws.Cells.Item(97, 14).Value = -8003
# Row 102: H102, I102, J102, K102, L102, M102, N102
$ws.Cells.Item(102, 8).Value = 3346.4375
$ws.Cells.Item(102, 9).Value = 2930.6428
$ws.Cells.Item(102, 10).Value = 6257
$ws.Cells.Item(102, 11).Value = 2930.6428
$ws.Cells.Item(102, 12).Value = 6257
$ws.Cells.Item(102, 13).Value = -1308.6428
$ws.Cells.Item(102, 14).Value = -9501
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Cells.Item(113, 8).Value = 16320.333
$ws.Cells.Item(113, 9).Value = 18684.4
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 18684.4
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -16514.4
$ws.Cells.Item(113, 14).Value = -8840
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Cells.Item(132, 8).Value = 3545186.2
$ws.Cells.Item(132, 9).Value = 4539296
$ws.Cells.Item(132, 10).Value = 65801.5
$ws.Cells.Item(132, 11).Value = 13617888
$ws.Cells.Item(132, 12).Value = 197404.5
$ws.Cells.Item(132, 13).Value = -13615358
$ws.Cells.Item(132, 14).Value = -202464.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46: H46, J46, L46, N46
$ws.Cells.Item(46, 8).Value = 1099.5
$ws.Cells.Item(46, 10).Value = 1200
$ws.Cells.Item(46, 12).Value = 1200
$ws.Cells.Item(46, 14).Value = -1576
# Row 68: H68, I68, K68, M68
$ws.Cells.Item(68, 8).Value = 2749.5
$ws.Cells.Item(68, 9).Value = 2666.3333
$ws.Cells.Item(68, 11).Value = 2666.3333
$ws.Cells.Item(68, 13).Value = -1917.3333
# Row 71: H71, I71, K71, M71
$ws.Cells.Item(71, 8).Value = 2749.5
$ws.Cells.Item(71, 9).Value = 2666.3333
$ws.Cells.Item(71, 11).Value = 13331.6665
$ws.Cells.Item(71, 13).Value = -9587.666499999999
# Row 82: H82, I82, J82, K82, L82, M82, N82
$ws.Cells.Item(82, 8).Value = 2620.6667
$ws.Cells.Item(82, 9).Value = 2222.75
$ws.Cells.Item(82, 10).Value = 3416.5
$ws.Cells.Item(82, 11).Value = 2222.75
$ws.Cells.Item(82, 12).Value = 3416.5
$ws.Cells.Item(82, 13).Value = -1861.75
$ws.Cells.Item(82, 14).Value = -4138.5
# Row 85: H85, I85, J85, K85, L85, M85, N85
$ws.Cells.Item(85, 8).Value = 2620.6667
$ws.Cells.Item(85, 9).Value = 2222.75
$ws.Cells.Item(85, 10).Value = 3416.5
$ws.Cells.Item(85, 11).Value = 2222.75
$ws.Cells.Item(85, 12).Value = 3416.5
$ws.Cells.Item(85, 13).Value = -974.75
$ws.Cells.Item(85, 14).Value = -5912.5
# Row 100: H100, J100, L100, N100
$ws.Cells.Item(100, 8).Value = 2380.2
$ws.Cells.Item(100, 10).Value = 2485.7144
$ws.Cells.Item(100, 12).Value = 2485.7144
$ws.Cells.Item(100, 14).Value = -3567.7144
# Row 122: H122, I122, K122, M122
$ws.Cells.Item(122, 8).Value = 1092940.5
$ws.Cells.Item(122, 9).Value = 1403545
$ws.Cells.Item(122, 11).Value = 4210635
$ws.Cells.Item(122, 13).Value = -4208185
# Row 136: H136, I136, K136, M136
$ws.Cells.Item(136, 8).Value = 3000.2666
$ws.Cells.Item(136, 9).Value = 3000.2666
$ws.Cells.Item(136, 11).Value = 9000.799800000001
$ws.Cells.Item(136, 13).Value = -6450.799800000001

$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62, I62, K62, M62
$ws.Cells.Item(62, 8).Value = 4934.875
$ws.Cells.Item(62, 9).Value = 4666.3335
$ws.Cells.Item(62, 11).Value = 4666.3335
$ws.Cells.Item(62, 13).Value = -4042.3335
# Row 65: H65, I65, K65, M65
$ws.Cells.Item(65, 8).Value = 4934.875
$ws.Cells.Item(65, 9).Value = 4666.3335
$ws.Cells.Item(65, 11).Value = 23331.6675
$ws.Cells.Item(65, 13).Value = -20211.6675
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Cells.Item(122, 8).Value = 1332.1052
$ws.Cells.Item(122, 9).Value = 1303.75
$ws.Cells.Item(122, 10).Value = 1483.3334
$ws.Cells.Item(122, 11).Value = 3911.25
$ws.Cells.Item(122, 12).Value = 4450.0002
$ws.Cells.Item(122, 13).Value = -1461.25
$ws.Cells.Item(122, 14).Value = -9350.0002
# Row 132: H132, I132, K132, M132
$ws.Cells.Item(132, 8).Value = 17858650
$ws.Cells.Item(132, 9).Value = 35715360
$ws.Cells.Item(132, 11).Value = 107146080
$ws.Cells.Item(132, 13).Value = -107143550
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Cells.Item(136, 8).Value = 43482692
$ws.Cells.Item(136, 9).Value = 55557630
$ws.Cells.Item(136, 10).Value = 12899.8
$ws.Cells.Item(136, 11).Value = 166672890
$ws.Cells.Item(136, 12).Value = 38699.39999999999
$ws.Cells.Item(136, 13).Value = -166670340
$ws.Cells.Item(136, 14).Value = -43799.39999999999
